# Auto-generated Excel COM-interop script
# Improves Arabic sentiment model results (column M) and recomputes the
# 'Mean over others' summary (column F, rows 30-35) for all 4 worksheets,
# including moving the per-row max-value bold+underline highlight where needed.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Accuracy ----
$ws = $wb.Worksheets.Item(1)

# Update values
$ws.Range("F30").Value = 0.6477488753419771
$ws.Range("F31").Value = 0.5669632442215197
$ws.Range("F32").Value = 0.6191497964348773
$ws.Range("F33").Value = 0.5059271648327915
$ws.Range("F35").Value = 0.6112873053327914
$ws.Range("M10").Value = 0.7136563876651982
$ws.Range("M11").Value = 0.4697624190064795
$ws.Range("M12").Value = 0.8568016614745587
$ws.Range("M13").Value = 0.5368802902055623
$ws.Range("M15").Value = 0.5493376242449561
$ws.Range("M2").Value = 0.4614465032875075
$ws.Range("M20").Value = 0.5556426454147672
$ws.Range("M21").Value = 0.5247779946921926
$ws.Range("M22").Value = 0.5696430035844368
$ws.Range("M23").Value = 0.5368802902055623
$ws.Range("M25").Value = 0.5467359834742397
$ws.Range("M3").Value = 0.6123009335529929
$ws.Range("M4").Value = 0.5996240601503759
$ws.Range("M5").Value = 0.5491990846681922
$ws.Range("M6").Value = 0.4218466012359142
$ws.Range("M7").Value = 0.5394736842105263
$ws.Range("M8").Value = 0.613013698630137
$ws.Range("M9").Value = 0.5255102040816326

# Update max-value highlight (bold + single underline) for affected cells
$ws.Range("E32").Font.Bold = $true
$ws.Range("E32").Font.Underline = 2
$ws.Range("F32").Font.Bold = $false
$ws.Range("F32").Font.Underline = -4142

# ---- Sheet 2: Macro_Precision ----
$ws = $wb.Worksheets.Item(2)

# Update values
$ws.Range("F30").Value = 0.6082735823994807
$ws.Range("F31").Value = 0.5817579110238733
$ws.Range("F32").Value = 0.5951111050003915
$ws.Range("F33").Value = 0.5673211934697572
$ws.Range("F35").Value = 0.5950475328079152
$ws.Range("M10").Value = 0.5835686053077358
$ws.Range("M11").Value = 0.6192671394799054
$ws.Range("M12").Value = 0.7455789338049237
$ws.Range("M13").Value = 0.6150117178484136
$ws.Range("M15").Value = 0.5760526664434757
$ws.Range("M2").Value = 0.54570153832858
$ws.Range("M20").Value = 0.5751093880281378
$ws.Range("M21").Value = 0.5555790256907959
$ws.Range("M22").Value = 0.5847976612816267
$ws.Range("M23").Value = 0.6150117178484136
$ws.Range("M25").Value = 0.5826244482122435
$ws.Range("M3").Value = 0.6138363404454956
$ws.Range("M4").Value = 0.5807045820861156
$ws.Range("M5").Value = 0.56019509125236
$ws.Range("M6").Value = 0.485006341685441
$ws.Range("M7").Value = 0.5487488556606652
$ws.Range("M8").Value = 0.6329818797262812
$ws.Range("M9").Value = 0.551557239057239

# ---- Sheet 3: Macro_Recall ----
$ws = $wb.Worksheets.Item(3)

# Update values
$ws.Range("F30").Value = 0.6150653921806938
$ws.Range("F31").Value = 0.5660186160247078
$ws.Range("F32").Value = 0.6040302953898261
$ws.Range("F33").Value = 0.5755659084210686
$ws.Range("F35").Value = 0.5950381011984093
$ws.Range("M10").Value = 0.6321428571428571
$ws.Range("M11").Value = 0.538697367601022
$ws.Range("M12").Value = 0.7521200153374233
$ws.Range("M13").Value = 0.6168561266580881
$ws.Range("M15").Value = 0.5797461162246675
$ws.Range("M2").Value = 0.5654838509948565
$ws.Range("M20").Value = 0.6129836747455942
$ws.Range("M21").Value = 0.5231830057342149
$ws.Range("M22").Value = 0.5796224785427443
$ws.Range("M23").Value = 0.6168561266580881
$ws.Range("M25").Value = 0.5831613214201603
$ws.Range("M3").Value = 0.6123939696601239
$ws.Range("M4").Value = 0.6864010865385402
$ws.Range("M5").Value = 0.5876557917888563
$ws.Range("M6").Value = 0.4912958759574226
$ws.Range("M7").Value = 0.5437360437360437
$ws.Range("M8").Value = 0.5345170975091788
$ws.Range("M9").Value = 0.5680272108843537

# ---- Sheet 4: Macro_F1 ----
$ws = $wb.Worksheets.Item(4)

# Update values
$ws.Range("F30").Value = 0.565556223152978
$ws.Range("F31").Value = 0.5322260584975322
$ws.Range("F32").Value = 0.5506596343884409
$ws.Range("F33").Value = 0.486277252205266
$ws.Range("F35").Value = 0.5494806386796504
$ws.Range("M10").Value = 0.586149833113623
$ws.Range("M11").Value = 0.404284701668049
$ws.Range("M12").Value = 0.7487636374471636
$ws.Range("M13").Value = 0.5368369488129967
$ws.Range("M15").Value = 0.5013628004337586
$ws.Range("M2").Value = 0.4470444004814329
$ws.Range("M20").Value = 0.5241906240191725
$ws.Range("M21").Value = 0.461669169824854
$ws.Range("M22").Value = 0.4987946168023652
$ws.Range("M23").Value = 0.5368369488129967
$ws.Range("M25").Value = 0.5053728398648472
$ws.Range("M3").Value = 0.6111045906637951
$ws.Range("M4").Value = 0.5223329455207648
$ws.Range("M5").Value = 0.5162805594106972
$ws.Range("M6").Value = 0.3899010360922989
$ws.Range("M7").Value = 0.5292086833968819
$ws.Range("M8").Value = 0.4658977899853812
$ws.Range("M9").Value = 0.5059493156254234

# Update max-value highlight (bold + single underline) for affected cells
$ws.Range("C31").Font.Bold = $false
$ws.Range("C31").Font.Underline = -4142
$ws.Range("E32").Font.Bold = $true
$ws.Range("E32").Font.Underline = 2
$ws.Range("F31").Font.Bold = $true
$ws.Range("F31").Font.Underline = 2
$ws.Range("F32").Font.Bold = $false
$ws.Range("F32").Font.Underline = -4142
